{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst p7 = paras.items[7];\nlet pf = p7.paragraphFormat;\nconst props = Object.getOwnPropertyNames(Object.getPrototypeOf(pf));\nreturn JSON.stringify(props);\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(9)\n$before = $p.Format.TabStops.Count\n$p.Format.TabStops.Add(999)\n$after = $p.Format.TabStops.Count\nWrite-Output \"before=$before after=$after\"\n"}
